$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column F
$ws.Range("F1").Value = "CRIT_RATE"

# Add new values for column F (rows 2-10)
$values = @(2, 3, 3, 3, 4, 4, 5, 5, 6)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}

# Update the selection as in the edited workbook
$ws.Range("G2").Select()
